$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table cell "A10[13]" -> the superscript footnote marker is split across
#    two runs ("[13" and "]") with a leftover "_GoBack" bookmark sitting in
#    between them. Re-typing the marker merges it back into a single run
#    "[13]" and drops the stray bookmark.
# ---------------------------------------------------------------------------
$scan = $d.Content
$scan.Start = 0
$scan.End = $d.Content.End
$targetStart = -1
$targetEnd = -1
while ($scan.Find.Execute("A10", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $e = $scan.End
    $peek = $d.Range($e, $e + 4)
    if ($peek.Text -eq "[13]") {
        $targetStart = $e
        $targetEnd = $e + 4
        break
    }
    $scan.Start = $e
    $scan.End = $d.Content.End
}
if ($targetStart -ge 0) {
    $rep = $d.Range($targetStart, $targetEnd)
    $rep.Find.Execute("[13]", $false, $false, $false, $false, $false, $true, 0, $false, "[13]", 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Table cell "A9 [18][19]" -> the " [18][19]" footnote markers are split
#    across three runs ("[", "18][", "19]") wrapped in a pair of gramStart/
#    gramEnd proof-error markers. Re-typing it merges everything back into a
#    single run and drops the proofErr markers.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(" [18][19]", $false, $false, $false, $false, $false, $true, 1, $false, " [18][19]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the prior sponsor's logo picture anchored in the page header.
# ---------------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($hi = 1; $hi -le $section.Headers.Count; $hi++) {
        $hdr = $section.Headers.Item($hi)
        for ($i = $hdr.Shapes.Count; $i -ge 1; $i--) {
            $hdr.Shapes.Item($i).Delete()
        }
    }
}
